$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of "cotações" (quotes) for 2025-10-29, mirroring the
# style/format of the preceding data row (row 54).
$newRow = 55

$ws.Cells.Item($newRow, 1).Value = Get-Date -Year 2025 -Month 10 -Day 29 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

$ws.Cells.Item($newRow, 2).Value = "22,0326"
$ws.Cells.Item($newRow, 3).Value = "15,9598"
$ws.Cells.Item($newRow, 4).Value = "15,9598"
$ws.Cells.Item($newRow, 5).Value = "15,9598"
